$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# C2: change from numeric 45 to the text "speed"
$ws.Range("C2").Value = "speed"

# B4: change from numeric 1 to the text "bird"
$ws.Range("B4").Value = "bird"

# C4: change from numeric 123 to numeric 123.34
$ws.Range("C4").Value = 123.34

# Row 5 (A5:C5) is removed entirely
$ws.Range("A5:C5").ClearContents()

# Update selection to match the post-edit state (C5, single cell)
$ws.Range("C5").Select()
